$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 43
$ws.Range("G2").Value = 67
$ws.Range("H2").Value = 48
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = $null

# Row 3
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = $null
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = $null
$ws.Range("K3").Value = $null

# Row 4
$ws.Range("F4").Value = 57
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 58
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = $null

# Row 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 27
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 38
$ws.Range("K5").Value = 50

# Row 6
$ws.Range("F6").Value = 39
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 34
$ws.Range("I6").Value = $null
$ws.Range("J6").Value = $null
$ws.Range("K6").Value = $null

# Row 7
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 22
$ws.Range("H7").Value = 40
$ws.Range("I7").Value = 60
$ws.Range("J7").Value = 64
$ws.Range("K7").Value = $null

# Row 8
$ws.Range("F8").Value = 29
$ws.Range("G8").Value = 65
$ws.Range("H8").Value = 36
$ws.Range("I8").Value = $null
$ws.Range("J8").Value = $null
$ws.Range("K8").Value = $null

# Row 9
$ws.Range("F9").Value = 33
$ws.Range("G9").Value = 18
$ws.Range("H9").Value = 40
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = $null
$ws.Range("K9").Value = $null

# Row 10
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 9
$ws.Range("H10").Value = 23
$ws.Range("I10").Value = 16
$ws.Range("J10").Value = 48
$ws.Range("K10").Value = 62

# Row 11
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 17
$ws.Range("J11").Value = 61
$ws.Range("K11").Value = 28

# Row 12
$ws.Range("F12").Value = 19
$ws.Range("G12").Value = 45
$ws.Range("H12").Value = 66
$ws.Range("I12").Value = $null
$ws.Range("J12").Value = $null
$ws.Range("K12").Value = $null

# Row 13
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = 28
$ws.Range("I13").Value = $null
$ws.Range("J13").Value = $null
$ws.Range("K13").Value = $null

# Row 14
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 46
$ws.Range("I14").Value = $null
$ws.Range("J14").Value = $null
$ws.Range("K14").Value = $null

# Row 15
$ws.Range("F15").Value = 63
$ws.Range("G15").Value = 38
$ws.Range("H15").Value = $null
$ws.Range("I15").Value = $null
$ws.Range("J15").Value = $null
$ws.Range("K15").Value = $null

# Row 16
$ws.Range("F16").Value = 57
$ws.Range("G16").Value = 56
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = $null
$ws.Range("K16").Value = $null

# Row 17
$ws.Range("F17").Value = 52
$ws.Range("G17").Value = 64
$ws.Range("H17").Value = $null
$ws.Range("I17").Value = $null
$ws.Range("J17").Value = $null
$ws.Range("K17").Value = $null

# Row 18
$ws.Range("F18").Value = 65
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 42
$ws.Range("I18").Value = $null
$ws.Range("J18").Value = $null
$ws.Range("K18").Value = $null

# Row 19
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = $null
$ws.Range("I19").Value = $null
$ws.Range("J19").Value = $null
$ws.Range("K19").Value = $null

# Row 20
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = $null
$ws.Range("I20").Value = $null
$ws.Range("J20").Value = $null
$ws.Range("K20").Value = $null

# Row 21
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = 10
$ws.Range("I21").Value = $null
$ws.Range("J21").Value = $null
$ws.Range("K21").Value = $null

# Row 22
$ws.Range("F22").Value = 37
$ws.Range("G22").Value = 53
$ws.Range("H22").Value = 59
$ws.Range("I22").Value = $null
$ws.Range("J22").Value = $null
$ws.Range("K22").Value = $null

# Row 23
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = 49
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = $null

# Row 24
$ws.Range("F24").Value = 23
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = 43
$ws.Range("I24").Value = 59
$ws.Range("J24").Value = 62
$ws.Range("K24").Value = $null

# Row 25
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = 32
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = $null

# Row 26
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 47
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = $null

# Row 27
$ws.Range("F27").Value = 19
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = $null
$ws.Range("I27").Value = $null
$ws.Range("J27").Value = $null
$ws.Range("K27").Value = $null

# Row 28
$ws.Range("F28").Value = 39
$ws.Range("G28").Value = 54
$ws.Range("H28").Value = $null
$ws.Range("I28").Value = $null
$ws.Range("J28").Value = $null
$ws.Range("K28").Value = $null

# Row 29
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 66
$ws.Range("H29").Value = $null
$ws.Range("I29").Value = $null
$ws.Range("J29").Value = $null
$ws.Range("K29").Value = $null

# Row 30
$ws.Range("F30").Value = 27
$ws.Range("G30").Value = 55
$ws.Range("H30").Value = $null
$ws.Range("I30").Value = $null
$ws.Range("J30").Value = $null
$ws.Range("K30").Value = $null

# Row 31
$ws.Range("F31").Value = 47
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 52
$ws.Range("I31").Value = 56
$ws.Range("J31").Value = 58
$ws.Range("K31").Value = $null

# Row 32
$ws.Range("F32").Value = 25
$ws.Range("G32").Value = 51
$ws.Range("H32").Value = 14
$ws.Range("I32").Value = 26
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = $null

# Row 33
$ws.Range("F33").Value = 21
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = $null
$ws.Range("I33").Value = $null
$ws.Range("J33").Value = $null
$ws.Range("K33").Value = $null

# Row 34
$ws.Range("F34").Value = 30
$ws.Range("G34").Value = 68
$ws.Range("H34").Value = $null
$ws.Range("I34").Value = $null
$ws.Range("J34").Value = $null
$ws.Range("K34").Value = $null

# Row 35
$ws.Range("F35").Value = 5
$ws.Range("G35").Value = 17
$ws.Range("H35").Value = 61
$ws.Range("I35").Value = $null
$ws.Range("J35").Value = $null
$ws.Range("K35").Value = $null

# Row 36
$ws.Range("F36").Value = 67
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 30
$ws.Range("I36").Value = $null
$ws.Range("J36").Value = $null
$ws.Range("K36").Value = $null

# Row 37
$ws.Range("F37").Value = 29
$ws.Range("G37").Value = 55
$ws.Range("H37").Value = $null
$ws.Range("I37").Value = $null
$ws.Range("J37").Value = $null
$ws.Range("K37").Value = $null

# Row 38
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 9
$ws.Range("H38").Value = 31
$ws.Range("I38").Value = 10
$ws.Range("J38").Value = $null
$ws.Range("K38").Value = $null

# Row 39
$ws.Range("F39").Value = 31
$ws.Range("G39").Value = 63
$ws.Range("H39").Value = 50
$ws.Range("I39").Value = 68
$ws.Range("J39").Value = $null
$ws.Range("K39").Value = $null

# Row 40
$ws.Range("F40").Value = 2
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 18
$ws.Range("I40").Value = 46
$ws.Range("J40").Value = $null
$ws.Range("K40").Value = $null

# Row 41
$ws.Range("F41").Value = 45
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 24
$ws.Range("I41").Value = 44
$ws.Range("J41").Value = $null
$ws.Range("K41").Value = $null
